# Auto-generated edit script applying the scheduled-runner profit recalculation update.
# For each (sheet,row) touched by the commit, set columns H-N to the new values;
# cells that are newly introduced are written directly, and cells that are removed
# in the target are cleared so no stale <c> element remains.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 700
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 733.3333
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 2199.9999
$ws.Range("M103").Value = -914
$ws.Range("N103").Value = -3371.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1708.3
$ws.Range("I112").Value = 2361.4285
$ws.Range("J112").Value = 1356.6154
$ws.Range("K112").Value = 7084.2855
$ws.Range("L112").Value = 4069.8462
$ws.Range("M112").Value = -5976.2855
$ws.Range("N112").Value = -6285.8462

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2565864.2
$ws.Range("I132").Value = 2741116.8
$ws.Range("J132").Value = 7180
$ws.Range("K132").Value = 8223350.399999999
$ws.Range("L132").Value = 21540
$ws.Range("M132").Value = -8220820.399999999
$ws.Range("N132").Value = -26600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18304.037
$ws.Range("I32").Value = 15149.372
$ws.Range("K32").Value = 15149.372
$ws.Range("M32").Value = -14862.372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2417.6667
$ws.Range("I63").Value = 1101.2
$ws.Range("J63").Value = 9000
$ws.Range("K63").Value = 1101.2
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -415.2
$ws.Range("N63").Value = -10372

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2417.6667
$ws.Range("I66").Value = 1101.2
$ws.Range("J66").Value = 9000
$ws.Range("K66").Value = 5506
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -2074
$ws.Range("N66").Value = -51864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 28074
$ws.Range("J35").Value = 28074
$ws.Range("L35").Value = 28074
$ws.Range("N35").Value = -28694

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20085.25
$ws.Range("J82").Value = 29176.6
$ws.Range("L82").Value = 29176.6
$ws.Range("N82").Value = -29942.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 20085.25
$ws.Range("J85").Value = 29176.6
$ws.Range("L85").Value = 29176.6
$ws.Range("N85").Value = -31828.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 4785.8
$ws.Range("J41").Value = 4910
$ws.Range("L41").Value = 4910
$ws.Range("N41").Value = -5766

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 24099
$ws.Range("J51").Value = 24099
$ws.Range("L51").Value = 24099
$ws.Range("N51").Value = -25571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 21487
$ws.Range("J59").Value = 21487
$ws.Range("L59").Value = 21487
$ws.Range("N59").Value = -23777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 16801.5
$ws.Range("J60").Value = 24103
$ws.Range("L60").Value = 24103
$ws.Range("N60").Value = -25125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 24099
$ws.Range("J61").Value = 24099
$ws.Range("L61").Value = 24099
$ws.Range("N61").Value = -24795

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 48295
$ws.Range("J68").Value = 48295
$ws.Range("L68").Value = 48295
$ws.Range("N68").Value = -49793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 48295
$ws.Range("J71").Value = 48295
$ws.Range("L71").Value = 144885
$ws.Range("N71").Value = -152373

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 24372.625
$ws.Range("J74").Value = 24372.625
$ws.Range("L74").Value = 24372.625
$ws.Range("N74").Value = -26120.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 24372.625
$ws.Range("J77").Value = 24372.625
$ws.Range("L77").Value = 73117.875
$ws.Range("N77").Value = -81853.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3281.9333
$ws.Range("I86").Value = 1903.2
$ws.Range("J86").Value = 6039.4
$ws.Range("K86").Value = 1903.2
$ws.Range("L86").Value = 6039.4
$ws.Range("M86").Value = -780.2
$ws.Range("N86").Value = -8285.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3281.9333
$ws.Range("I89").Value = 1903.2
$ws.Range("J89").Value = 6039.4
$ws.Range("K89").Value = 9516
$ws.Range("L89").Value = 30197
$ws.Range("M89").Value = -3900
$ws.Range("N89").Value = -41429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1410.8197
$ws.Range("J131").Value = 1152.1698
$ws.Range("L131").Value = 3456.5094
$ws.Range("N131").Value = -13536.5094

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2482.6
$ws.Range("I43").Value = 1225.875
$ws.Range("J43").Value = 7509.5
$ws.Range("K43").Value = 1225.875
$ws.Range("L43").Value = 7509.5
$ws.Range("M43").Value = -1074.875
$ws.Range("N43").Value = -7811.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 17780.5
$ws.Range("J57").Value = 24061
$ws.Range("L57").Value = 24061
$ws.Range("N57").Value = -25701

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3267.7778
$ws.Range("I80").Value = 3325.3333
$ws.Range("J80").Value = 2980
$ws.Range("K80").Value = 3325.3333
$ws.Range("L80").Value = 2980
$ws.Range("M80").Value = -2327.3333
$ws.Range("N80").Value = -4976

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3267.7778
$ws.Range("I83").Value = 3325.3333
$ws.Range("J83").Value = 2980
$ws.Range("K83").Value = 16626.6665
$ws.Range("L83").Value = 14900
$ws.Range("M83").Value = -11634.6665
$ws.Range("N83").Value = -24884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9650.799999999999
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 9650.799999999999
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 28952.4
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -33852.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3384.0356
$ws.Range("I132").Value = 3349.5898
$ws.Range("K132").Value = 10048.7694
$ws.Range("M132").Value = -7518.769400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2912.3809
$ws.Range("I132").Value = 2209.5652
$ws.Range("J132").Value = 3763.158
$ws.Range("K132").Value = 6628.6956
$ws.Range("L132").Value = 11289.474
$ws.Range("M132").Value = -4098.6956
$ws.Range("N132").Value = -16349.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2946848
$ws.Range("I136").Value = 5268169.5
$ws.Range("J136").Value = 6507.3335
$ws.Range("K136").Value = 15804508.5
$ws.Range("L136").Value = 19522.0005
$ws.Range("M136").Value = -15801958.5
$ws.Range("N136").Value = -24622.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2110.9092
$ws.Range("I136").Value = 1625.5128
$ws.Range("J136").Value = 3294.0625
$ws.Range("K136").Value = 4876.538399999999
$ws.Range("L136").Value = 9882.1875
$ws.Range("M136").Value = -2326.538399999999
$ws.Range("N136").Value = -14982.1875

